$wb = $excel.ActiveWorkbook

# --- registration sheet: add rows 4-6 ---
$wsReg = $wb.Worksheets.Item("registration")

$wsReg.Range("A4").Value = "harh"
$wsReg.Range("B4").Value = "nayak"
$wsReg.Range("C4").Value = "harh1@example.com"
$wsReg.Range("D4").Value = "'9876543212"
$wsReg.Range("E4").Value = "harh123"
$wsReg.Range("F4").Value = "harh123"
$wsReg.Range("G4").Value = "Yes"

$wsReg.Range("A5").Value = "jane"
$wsReg.Range("B5").Value = "smith"
$wsReg.Range("C5").Value = "jane@example.com"
$wsReg.Range("D5").Value = "'9876543211"
$wsReg.Range("E5").Value = "jane123"
$wsReg.Range("F5").Value = "jane123"
$wsReg.Range("G5").Value = "Yes"

$wsReg.Range("A6").Value = "dave"
$wsReg.Range("B6").Value = "franc"
$wsReg.Range("C6").Value = "dave@example.com"
$wsReg.Range("D6").Value = "'9876543215"
$wsReg.Range("E6").Value = "dave123"
$wsReg.Range("F6").Value = "dave123"
$wsReg.Range("G6").Value = "Yes"

# --- login sheet: add rows 3-5 ---
$wsLogin = $wb.Worksheets.Item("login")

$wsLogin.Range("A3").Value = "harh1@example.com"
$wsLogin.Range("B3").Value = "harh123"
$wsLogin.Range("C3").Value = "Login Successful"

$wsLogin.Range("A4").Value = "anand1@example.com"
$wsLogin.Range("B4").Value = "anand123"
$wsLogin.Range("C4").Value = "Login Successful"

$wsLogin.Range("A5").Value = "pratham1@example.com"
$wsLogin.Range("B5").Value = "pratham123"
$wsLogin.Range("C5").Value = "Login Successful"

# --- search sheet: add rows 4-8 ---
$wsSearch = $wb.Worksheets.Item("search")

$wsSearch.Range("A4").Value = "samsung"
$wsSearch.Range("A5").Value = "sony"
$wsSearch.Range("A6").Value = "canon"
$wsSearch.Range("A7").Value = "samsung"
$wsSearch.Range("A8").Value = "tab"
